$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add 6 new locator rows (23-28) to the ContactPage sheet ---

# Row 23: ContactPage_TextBox_Zip
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = "input#PostalCode"
$ws.Range("B23").Value = "ContactPage_TextBox_Zip"

# Row 24: ContactPage_ErrorMessage_InvalidZip_TextBox_Zip (text-formatted locator name)
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "ContactPage_ErrorMessage_InvalidZip_TextBox_Zip"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("C24").Value = ".mktoError #ValidMsgPostalCode"

# Row 25: ContactPage_Dropdown_Province
$ws.Range("A25").Value = 24
$ws.Range("C25").Value = "select#State"
$ws.Range("B25").Value = "ContactPage_Dropdown_Province"

# Row 26: ContactPage_ErrorMessage_InvalidProvince_Dropdown_Province (text-formatted locator name)
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "ContactPage_ErrorMessage_InvalidProvince_Dropdown_Province"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("C26").Value = ".mktoError #ValidMsgState"

# Row 27: ContactPage_CheckBox_PrivacyConsent (text-formatted locator name)
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "ContactPage_CheckBox_PrivacyConsent"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("C27").Value = "input#mktoCheckbox_142098_0"

# Row 28: ContactPage_Label_PrivacyConsent (text-formatted locator name)
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "ContactPage_Label_PrivacyConsent"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("C28").Value = "label#LblmktoCheckbox_142098_0"

# --- Update the selected range shown when the sheet is opened ---
$ws.Range("B13:C13").Select()

# --- Switch the workbook off manual calculation (back to automatic) ---
$excel.Calculation = -4105
